$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at the top of the "Poroto verde" data
# block (row 199), pushing all the existing rows (199-232) down by one
# (to 200-233). Insert a blank row at 199 to reproduce that shift.
$ws.Rows("199:199").Insert()

# Populate the newly inserted row 199 with the new record's data.
$ws.Range("A199").Value = 5
$ws.Range("B199").Value = "Macroferia Regional de Talca"
$ws.Range("C199").Value = "Maule"
$ws.Range("D199").Value = 44995
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 100112031
$ws.Range("G199").Value = "Poroto verde"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 150
$ws.Range("K199").Value = 25000
$ws.Range("L199").Value = 25000
$ws.Range("M199").Value = 25000
$ws.Range("N199").Value = "$/saco 25 kilos"
$ws.Range("O199").Value = "Región del Maule"
$ws.Range("P199").Value = 1000
$ws.Range("Q199").Value = 25
$ws.Range("R199").Value = "Hortaliza"
